$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Id (A), Antal (I), and Publik kommentar (AC) values between
# row 2 and row 3, and replace the Ost (Q) / Nord (R) coordinates with
# their rounded, swapped counterparts.
$ws.Range("A2").Value = 111676783
$ws.Range("A3").Value = 111676874

$ws.Range("I2").Value = "'100"
$ws.Range("I3").Value = "'15"

$ws.Range("Q2").Value = 580098
$ws.Range("R2").Value = 6414648
$ws.Range("Q3").Value = 580088
$ws.Range("R3").Value = 6414641

$ws.Range("AC2").Value = "3 blommor"
$ws.Range("AC3").Value = "1 blomma"

# Remove the now-empty Starttid (Z) / Sluttid (AB) cells for both rows.
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
